$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.542.10"
$ws.Range("E2").Value = "  +0.79%  "
$ws.Range("D3").Value = "3.619.14"
$ws.Range("E3").Value = "  +2.36%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = $ws.Range("B4").Style
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.95"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value = "  +0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "196.73"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E9").Value = "  -0.58%  "
$ws.Range("E10").Value = "  -0.91%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "53.63"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value = "  -0.48%  "
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "9.57"
$ws.Range("D13").Style = $ws.Range("B13").Style
$ws.Range("E13").Value = "  +0.46%  "
$ws.Range("D14").Value = "4.197.24"
$ws.Range("E14").Value = "  +2.21%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "602.26"
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Value = "  +0.07%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.00"
$ws.Range("D16").Style = $ws.Range("B16").Style
$ws.Range("E16").Value = "  +2.02%  "
$ws.Range("D17").Value = "70.706.90"
$ws.Range("E17").Value = "  +0.77%  "
$ws.Range("D18").Value = "3.622.04"
$ws.Range("E18").Value = "  +2.31%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "19.09"
$ws.Range("D19").Style = $ws.Range("B19").Style
$ws.Range("E19").Value = "  -0.46%  "
$ws.Range("E20").Value = "  +1.49%  "
$ws.Range("E21").Value = "  +0.57%  "
$ws.Range("E22").Value = "  -0.53%  "
$ws.Range("E23").Value = "  -1.29%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "101.95"
$ws.Range("D24").Style = $ws.Range("B24").Style
$ws.Range("E24").Value = "  -0.70%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "4.63"
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("E25").Value = "  +0.34%  "
$ws.Range("E26").Value = "  -3.82%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.77"
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value = "  -1.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.64"
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value = "  +0.47%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.89"
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value = "  +0.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.68"
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value = "  +7.25%  "
$ws.Range("E31").Value = "  +1.61%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "12.29"
$ws.Range("D32").Style = $ws.Range("B32").Style
$ws.Range("E32").Value = "  -2.49%  "
$ws.Range("E33").Value = "  +2.49%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "63.56"
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("D35").Value = "0.0₃0892"
$ws.Range("E35").Value = "  +7.14%  "
$ws.Range("D36").Value = "3.917.45"
$ws.Range("E36").Value = "  +3.65%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "539.72"
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value = "  +9.48%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.11"
$ws.Range("D38").Style = $ws.Range("B38").Style
$ws.Range("E38").Value = "  +1.59%  "
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("E40").Value = "  +1.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.390"
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value = "  -0.89%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.55"
$ws.Range("D42").Style = $ws.Range("B42").Style
$ws.Range("E42").Value = "  -1.87%  "
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0461"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value = "  +0.77%  "
$ws.Range("E45").Value = "  +3.24%  "
$ws.Range("E46").Value = "  +0.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.140"
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value = "  +0.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.62"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value = "  -0.33%  "
$ws.Range("E49").Value = "  -0.20%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000251"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value = "  +0.96%  "
$ws.Range("E51").Value = "  +3.07%  "
